$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 176
$ws.Range("B176").Value = 7302796
$ws.Range("F176").Value = "Sport Huancayo"
$ws.Range("G176").Value = "Sport Boys"
$ws.Range("H176").Value = 1
$ws.Range("I176").Value = 0
$ws.Range("J176").Value = "H"
$ws.Range("K176").Value = 1.727
$ws.Range("L176").Value = 3.75
$ws.Range("M176").Value = 4.333
$ws.Range("N176").Value = 1.25
$ws.Range("O176").Value = 5.25
$ws.Range("P176").Value = 10
$ws.Range("Q176").Value = -1.75
$ws.Range("R176").Value = 1.925
$ws.Range("S176").Value = 1.875
$ws.Range("T176").Value = 3
$ws.Range("U176").Value = 1.875
$ws.Range("V176").Value = 1.925
$ws.Range("W176").Value = 0.25
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = -1
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 0.875
$ws.Range("AB176").Value = -1
$ws.Range("AC176").Value = 0.925

# Row 177
$ws.Range("B177").Value = 7302795
$ws.Range("F177").Value = "Unin Comercio"
$ws.Range("G177").Value = "Deportivo Garcilaso"
$ws.Range("H177").Value = 1
$ws.Range("I177").Value = 2
$ws.Range("J177").Value = "A"
$ws.Range("K177").Value = 2.25
$ws.Range("L177").Value = 3.3
$ws.Range("M177").Value = 2.7
$ws.Range("N177").Value = 1.75
$ws.Range("O177").Value = 3.6
$ws.Range("P177").Value = 4
$ws.Range("Q177").Value = -0.5
$ws.Range("R177").Value = 1.8
$ws.Range("S177").Value = 2
$ws.Range("T177").Value = 2.75
$ws.Range("U177").Value = 1.825
$ws.Range("V177").Value = 1.975
$ws.Range("W177").Value = -1
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = 3
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 1
$ws.Range("AB177").Value = 0.4125
$ws.Range("AC177").Value = -0.5

# Row 180
$ws.Range("B180").Value = 7384624
$ws.Range("F180").Value = "Cesar Vallejo"
$ws.Range("G180").Value = "Cusco FC"
$ws.Range("H180").Value = 3
$ws.Range("I180").Value = 1
$ws.Range("J180").Value = "H"
$ws.Range("K180").Value = 2
$ws.Range("L180").Value = 3.4
$ws.Range("M180").Value = 3.5
$ws.Range("N180").Value = 1.45
$ws.Range("O180").Value = 4.2
$ws.Range("P180").Value = 6.5
$ws.Range("Q180").Value = -1
$ws.Range("R180").Value = 1.75
$ws.Range("S180").Value = 2.05
$ws.Range("T180").Value = 2.5
$ws.Range("U180").Value = 1.95
$ws.Range("V180").Value = 1.85
$ws.Range("W180").Value = 0.45
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.75
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = 0.95
$ws.Range("AC180").Value = -1

# Row 181
$ws.Range("B181").Value = 7384623
$ws.Range("F181").Value = "Sport Boys"
$ws.Range("G181").Value = "Cienciano"
$ws.Range("H181").Value = 2
$ws.Range("I181").Value = 1
$ws.Range("J181").Value = "H"
$ws.Range("K181").Value = 2
$ws.Range("L181").Value = 3.4
$ws.Range("M181").Value = 3.5
$ws.Range("N181").Value = 1.833
$ws.Range("O181").Value = 4
$ws.Range("P181").Value = 3.2
$ws.Range("Q181").Value = -0.5
$ws.Range("R181").Value = 1.925
$ws.Range("S181").Value = 1.875
$ws.Range("T181").Value = 3
$ws.Range("U181").Value = 1.925
$ws.Range("V181").Value = 1.875
$ws.Range("W181").Value = 0.833
$ws.Range("X181").Value = -1
$ws.Range("Y181").Value = -1
$ws.Range("Z181").Value = 0.925
$ws.Range("AA181").Value = -1
$ws.Range("AB181").Value = 0
$ws.Range("AC181").Value = -0

# Row 182
$ws.Range("B182").Value = 7384622
$ws.Range("F182").Value = "Deportivo Municipal"
$ws.Range("G182").Value = "Academia Deportiva Cantolao"
$ws.Range("H182").Value = 1
$ws.Range("I182").Value = 2
$ws.Range("J182").Value = "A"
$ws.Range("K182").Value = 1.444
$ws.Range("L182").Value = 4.333
$ws.Range("M182").Value = 7
$ws.Range("N182").Value = 1.5
$ws.Range("O182").Value = 3.75
$ws.Range("P182").Value = 6
$ws.Range("Q182").Value = -1
$ws.Range("R182").Value = 1.825
$ws.Range("S182").Value = 2.025
$ws.Range("T182").Value = 2.75
$ws.Range("U182").Value = 1.875
$ws.Range("V182").Value = 1.975
$ws.Range("W182").Value = -1
$ws.Range("X182").Value = -1
$ws.Range("Y182").Value = 5
$ws.Range("Z182").Value = -1
$ws.Range("AA182").Value = 1.025
$ws.Range("AB182").Value = 0.4375
$ws.Range("AC182").Value = -0.5

# Row 183
$ws.Range("B183").Value = 7384630
$ws.Range("F183").Value = "Atletico Grau"
$ws.Range("G183").Value = "Unin Comercio"
$ws.Range("H183").Value = 0
$ws.Range("I183").Value = 1
$ws.Range("J183").Value = "A"
$ws.Range("K183").Value = 2.8
$ws.Range("L183").Value = 3.4
$ws.Range("M183").Value = 2.15
$ws.Range("N183").Value = 1.75
$ws.Range("O183").Value = 3.6
$ws.Range("P183").Value = 3.8
$ws.Range("Q183").Value = -0.75
$ws.Range("R183").Value = 2
$ws.Range("S183").Value = 1.8
$ws.Range("T183").Value = 3
$ws.Range("U183").Value = 1.85
$ws.Range("V183").Value = 1.95
$ws.Range("W183").Value = -1
$ws.Range("X183").Value = -1
$ws.Range("Y183").Value = 2.8
$ws.Range("Z183").Value = -1
$ws.Range("AA183").Value = 0.8
$ws.Range("AB183").Value = -1
$ws.Range("AC183").Value = 0.95

# Row 184
$ws.Range("B184").Value = 7384627
$ws.Range("F184").Value = "Universitario de Deportes"
$ws.Range("G184").Value = "Sport Huancayo"
$ws.Range("H184").Value = 2
$ws.Range("I184").Value = 0
$ws.Range("J184").Value = "H"
$ws.Range("K184").Value = 1.25
$ws.Range("L184").Value = 5
$ws.Range("M184").Value = 12
$ws.Range("N184").Value = 1.181
$ws.Range("O184").Value = 6
$ws.Range("P184").Value = 13
$ws.Range("Q184").Value = -1.75
$ws.Range("R184").Value = 1.8
$ws.Range("S184").Value = 2
$ws.Range("T184").Value = 2.75
$ws.Range("U184").Value = 1.85
$ws.Range("V184").Value = 1.95
$ws.Range("W184").Value = 0.181
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = 0.4
$ws.Range("AA184").Value = -0.5
$ws.Range("AB184").Value = -1
$ws.Range("AC184").Value = 0.95

# Row 185
$ws.Range("B185").Value = 7384629
$ws.Range("F185").Value = "Deportivo Garcilaso"
$ws.Range("G185").Value = "Alianza Lima"
$ws.Range("H185").Value = 0
$ws.Range("I185").Value = 1
$ws.Range("J185").Value = "A"
$ws.Range("K185").Value = 2.625
$ws.Range("L185").Value = 3.3
$ws.Range("M185").Value = 2.5
$ws.Range("N185").Value = 2.7
$ws.Range("O185").Value = 3.4
$ws.Range("P185").Value = 2.375
$ws.Range("Q185").Value = 0
$ws.Range("R185").Value = 2.025
$ws.Range("S185").Value = 1.775
$ws.Range("T185").Value = 2.25
$ws.Range("U185").Value = 1.825
$ws.Range("V185").Value = 1.975
$ws.Range("W185").Value = -1
$ws.Range("X185").Value = -1
$ws.Range("Y185").Value = 1.375
$ws.Range("Z185").Value = -1
$ws.Range("AA185").Value = 0.7749999999999999
$ws.Range("AB185").Value = -1
$ws.Range("AC185").Value = 0.9750000000000001

# Row 186
$ws.Range("B186").Value = 7384625
$ws.Range("F186").Value = "AD Tarma"
$ws.Range("G186").Value = "Carlos Manucci"
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = 0
$ws.Range("J186").Value = "D"
$ws.Range("K186").Value = 1.5
$ws.Range("L186").Value = 3.75
$ws.Range("M186").Value = 7
$ws.Range("N186").Value = 1.363
$ws.Range("O186").Value = 4.333
$ws.Range("P186").Value = 9.5
$ws.Range("Q186").Value = -1.25
$ws.Range("R186").Value = 1.875
$ws.Range("S186").Value = 1.925
$ws.Range("T186").Value = 2.5
$ws.Range("U186").Value = 1.8
$ws.Range("V186").Value = 2
$ws.Range("W186").Value = -1
$ws.Range("X186").Value = 3.333
$ws.Range("Y186").Value = -1
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = 0.925
$ws.Range("AB186").Value = -1
$ws.Range("AC186").Value = 1

# Row 187
$ws.Range("B187").Value = 7384628
$ws.Range("F187").Value = "Deportivo Binacional"
$ws.Range("G187").Value = "FBC Melgar"
$ws.Range("H187").Value = 1
$ws.Range("I187").Value = 2
$ws.Range("J187").Value = "A"
$ws.Range("K187").Value = 2.75
$ws.Range("L187").Value = 3.3
$ws.Range("M187").Value = 2.375
$ws.Range("N187").Value = 3.3
$ws.Range("O187").Value = 3.6
$ws.Range("P187").Value = 2
$ws.Range("Q187").Value = 0.5
$ws.Range("R187").Value = 1.8
$ws.Range("S187").Value = 2
$ws.Range("T187").Value = 2.75
$ws.Range("U187").Value = 1.975
$ws.Range("V187").Value = 1.875
$ws.Range("W187").Value = -1
$ws.Range("X187").Value = -1
$ws.Range("Y187").Value = 1
$ws.Range("Z187").Value = -1
$ws.Range("AA187").Value = 1
$ws.Range("AB187").Value = 0.4875
$ws.Range("AC187").Value = -0.5

# Row 188
$ws.Range("B188").Value = 7384626
$ws.Range("F188").Value = "Sporting Cristal"
$ws.Range("G188").Value = "Alianza Atletico"
$ws.Range("H188").Value = 3
$ws.Range("I188").Value = 0
$ws.Range("J188").Value = "H"
$ws.Range("K188").Value = 1.3
$ws.Range("L188").Value = 5
$ws.Range("M188").Value = 9
$ws.Range("N188").Value = 1.166
$ws.Range("O188").Value = 6.5
$ws.Range("P188").Value = 13
$ws.Range("Q188").Value = -2
$ws.Range("R188").Value = 1.85
$ws.Range("S188").Value = 1.95
$ws.Range("T188").Value = 3.25
$ws.Range("U188").Value = 2
$ws.Range("V188").Value = 1.8
$ws.Range("W188").Value = 0.1659999999999999
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = -1
$ws.Range("Z188").Value = 0.8500000000000001
$ws.Range("AA188").Value = -1
$ws.Range("AB188").Value = -0.5
$ws.Range("AC188").Value = 0.4

# Row 228
$ws.Range("B228").Value = 7818816
$ws.Range("F228").Value = "UTC Cajamarca"
$ws.Range("G228").Value = "Universitario de Deportes"
$ws.Range("H228").Value = 0
$ws.Range("I228").Value = 0
$ws.Range("J228").Value = "D"
$ws.Range("K228").Value = 3.3
$ws.Range("L228").Value = 3.3
$ws.Range("M228").Value = 2.1
$ws.Range("N228").Value = 4.5
$ws.Range("O228").Value = 3.2
$ws.Range("P228").Value = 1.95
$ws.Range("Q228").Value = 0.5
$ws.Range("R228").Value = 2
$ws.Range("S228").Value = 1.85
$ws.Range("T228").Value = 2
$ws.Range("U228").Value = 1.775
$ws.Range("V228").Value = 2.1
$ws.Range("W228").Value = -1
$ws.Range("X228").Value = 2.2
$ws.Range("Y228").Value = -1
$ws.Range("Z228").Value = 1
$ws.Range("AA228").Value = -1
$ws.Range("AB228").Value = -1
$ws.Range("AC228").Value = 1.1

# Row 229
$ws.Range("B229").Value = 7818817
$ws.Range("F229").Value = "Sport Boys"
$ws.Range("G229").Value = "Cusco FC"
$ws.Range("H229").Value = 3
$ws.Range("I229").Value = 0
$ws.Range("J229").Value = "H"
$ws.Range("K229").Value = 2.2
$ws.Range("L229").Value = 3.2
$ws.Range("M229").Value = 3.2
$ws.Range("N229").Value = 1.6
$ws.Range("O229").Value = 3.75
$ws.Range("P229").Value = 5.75
$ws.Range("Q229").Value = -0.75
$ws.Range("R229").Value = 1.85
$ws.Range("S229").Value = 2
$ws.Range("T229").Value = 2.5
$ws.Range("U229").Value = 1.975
$ws.Range("V229").Value = 1.875
$ws.Range("W229").Value = 0.6000000000000001
$ws.Range("X229").Value = -1
$ws.Range("Y229").Value = -1
$ws.Range("Z229").Value = 0.8500000000000001
$ws.Range("AA229").Value = -1
$ws.Range("AB229").Value = 0.9750000000000001
$ws.Range("AC229").Value = -1

# Row 280
$ws.Range("H280").Value = 0
$ws.Range("I280").Value = 1
$ws.Range("J280").Value = "A"
$ws.Range("N280").Value = 2.4
$ws.Range("O280").Value = 3.1
$ws.Range("P280").Value = 3
$ws.Range("R280").Value = 2.05
$ws.Range("S280").Value = 1.75
$ws.Range("V280").Value = 1.85
$ws.Range("W280").Value = -1
$ws.Range("X280").Value = -1
$ws.Range("Y280").Value = 2
$ws.Range("Z280").Value = -1
$ws.Range("AA280").Value = 0.75
$ws.Range("AB280").Value = -1
$ws.Range("AC280").Value = 0.8500000000000001
